# Fix incorrect dates / identity in the generated invoice (generative case),
# and fix the local invoice address so it renders correctly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block: name / dates / contract / invoice number ---
$ws.Range("A1").Value = "Ivanov Ivan Invoice"
$ws.Range("A2").Value = "Invoice date: October 1, 2000"
$ws.Range("A3").Value = "Contract: dated as of September 1, 2000"
$ws.Range("A4").Value = "Invoice number: 2000-10-II"
$ws.Range("A5").Value = "Date of service: October 2000"

# --- Bank address / contact line ---
$ws.Range("B21").Value = "1 Lenina str., Moscow, 1000000, tel +7 495 755-58-58, SWIFT "

# --- Beneficiary address: fix the postal code, keeping the zip as its own run ---
$ws.Range("B29").Value = "PR. LENINA, D. 1, KV. 1, MOSCOW, RUSSIA, 1000000"
$zip = $ws.Range("B29").Characters(42, 7)
$zip.Font.Name = "Arial"
$zip.Font.Size = 10

# --- Selection restored to where the user was working ---
$ws.Range("A21").Select() | Out-Null
